$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-10 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-11 Wednesday", 2)
$d.Content.Find.Execute("702÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "857÷6=", 2)
$d.Content.Find.Execute("701÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "328÷3=", 2)
$d.Content.Find.Execute("645÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷9=", 2)
$d.Content.Find.Execute("513÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "995÷9=", 2)
$d.Content.Find.Execute("924÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "880÷2=", 2)
$d.Content.Find.Execute("501÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "489÷5=", 2)
$d.Content.Find.Execute("822÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "892÷4=", 2)
$d.Content.Find.Execute("402÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "360÷2=", 2)
$d.Content.Find.Execute("501÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "552÷9=", 2)
$d.Content.Find.Execute("299÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "781÷9=", 2)
$d.Content.Find.Execute("222÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "398÷6=", 2)
$d.Content.Find.Execute("196÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "358÷2=", 2)
$d.Content.Find.Execute("797÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "963÷7=", 2)
$d.Content.Find.Execute("490÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "649÷7=", 2)
$d.Content.Find.Execute("636÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "484÷2=", 2)
$d.Content.Find.Execute("445÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "916÷7=", 2)
$d.Content.Find.Execute("637÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "377÷8=", 2)
$d.Content.Find.Execute("381÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "785÷7=", 2)
$d.Content.Find.Execute("249÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "547÷6=", 2)
$d.Content.Find.Execute("975÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "783÷2=", 2)
$d.Content.Find.Execute("115÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "883÷9=", 2)
$d.Content.Find.Execute("233÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "117÷3=", 2)
$d.Content.Find.Execute("353÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "701÷3=", 2)
$d.Content.Find.Execute("833÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "546÷8=", 2)
$d.Content.Find.Execute("516÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "637÷5=", 2)
